#
# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The "Rules" decision table's R40 rule-name cell (B11) is renamed to "1".
# The new literal text is stored (like the original "R40") as a shared
# string, so we force the cell to Text before writing the value - otherwise
# a bare "1" would be auto-recognised as a number.
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
